# Update the "as of" date in the confidential disclosure blurb (A18) and
# refresh the Weight / Percent Change figures for rows 2-15 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect so the cells can be written, then
# restore protection (same flags as before) once the edits are done.
$ws.Unprotect()

# --- Update the confidential disclosure text (A18): 2021-04-26 -> 2021-04-27 ---
$oldText = $ws.Range("A18").Value2
$newText = $oldText -replace '2021-04-26', '2021-04-27'
$ws.Range("A18").Value = $newText
# Undo the auto row-height adjustment that Excel applies when the wrapped
# text is rewritten, so row 18 goes back to the sheet's default height.
$ws.Rows("18").AutoFit()

# --- Refresh Weight (D) / Percent Change (E) values for rows 2-15 ---
$newValues = @{
    2  = @(0.05759422717737765,  -0.0004587155963302614)
    3  = @(0.02370874287567831,  -0.003120124804992352)
    4  = @(0.03173098465660066,   0.001498688647433521)
    5  = @(0.03050710732908138,   0.008227067050596615)
    6  = @(0.03726716369560042,  -0.004891535516801326)
    7  = @(0.01903028850685467,   0.007147559941595105)
    8  = @(0.004878336719404947, -0.003249390739236357)
    9  = @(0.006940368568576656, -0.002664636467453407)
    10 = @(0.07014342805318242,  -0.002824858757062176)
    11 = @(0.07026231521937425,  -0.002820078962210926)
    12 = @(0.1478005250096888,   -0.008723008723008907)
    13 = @(0.3853192499860456,   -0.00393391030684509)
    14 = @(0.1148172622025341,   -0.003382459531287796)
    15 = @(0.9999999999999999,   -0.003472205366113812)
}

foreach ($row in $newValues.Keys) {
    $pair = $newValues[$row]
    $ws.Range("D$row").Value = $pair[0]
    $ws.Range("E$row").Value = $pair[1]
}

# Restore sheet protection with the same settings that were in effect
# before the edit (contents + objects + scenarios protected).
$ws.Protect($null, $true, $true, $true)
